$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, shifting rows 159:256 down to 160:257.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with its data.
$ws.Range("A159").Value = 7
$ws.Range("B159").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C159").Value = "Ñuble"
$ws.Range("D159").Value = 44879
$ws.Range("E159").Value = 16
$ws.Range("F159").Value = 100112032
$ws.Range("G159").Value = "Zapallo italiano"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 120
$ws.Range("K159").Value = 9000
$ws.Range("L159").Value = 10000
$ws.Range("M159").Value = 9500
$ws.Range("N159").Value = "$/caja 50 unidades"
$ws.Range("O159").Value = "Región de O'Higgins"
$ws.Range("P159").Value = 190
$ws.Range("Q159").Value = 50
$ws.Range("R159").Value = "Hortaliza"

# Keep the date-formatted number format consistent with other rows in column D.
$ws.Range("D159").NumberFormat = "YYYY-MM-DD HH:MM:SS"
